# For every row in the used range of the active sheet, swap the values of
# column D (codeforiati:category-name) and column G (codeforiati:group-code).
# Columns E (codeforiati:group-name) and F (codeforiati:category-code) are
# left untouched. This matches the reordering applied to the "before" file
# where D and G were transposed across the whole SectorGroup table
# (including the header row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ur = $ws.UsedRange
$rowCount = $ur.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $gCell = $ws.Cells.Item($r, 7)

    $dVal = $dCell.Value()
    $gVal = $gCell.Value()

    $dCell.Value = $gVal
    $gCell.Value = $dVal
}
